$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(4, "Porphyromonas", 0.00929406095537628),
    @(5, "Haemophilus", 0.005704838249311408),
    @(6, "Neisseria", 0.0003247717504744928),
    @(8, "Porphyromonas", 0.00929406095537628),
    @(9, "Varibaculum", 0.006990870234154651),
    @(10, "Ezakiella", 0.002980975725860719),
    @(16, "Porphyromonas", 0.00929406095537628),
    @(17, "Haemophilus", 0.005704838249311408),
    @(18, "Neisseria", 0.0003247717504744928),
    @(20, "Porphyromonas", 0.00929406095537628),
    @(21, "Varibaculum", 0.006990870234154651),
    @(22, "Ezakiella", 0.002980975725860719),
    @(41, "Porphyromonas", 0.00929406095537628),
    @(42, "Neisseria", 0.0003247717504744928),
    @(45, "Porphyromonas", 0.00929406095537628),
    @(46, "Ezakiella", 0.002980975725860719),
    @(55, "Prevotella", 0.05315508228653677),
    @(56, "Porphyromonas", 0.00929406095537628),
    @(57, "Varibaculum", 0.006990870234154651),
    @(58, "Ezakiella", 0.002980975725860719),
    @(64, "Porphyromonas", 0.00929406095537628),
    @(65, "Haemophilus", 0.005704838249311408),
    @(66, "Neisseria", 0.0003247717504744928),
    @(68, "Porphyromonas", 0.00929406095537628),
    @(69, "Varibaculum", 0.006990870234154651),
    @(70, "Ezakiella", 0.002980975725860719),
    @(89, "Porphyromonas", 0.00929406095537628),
    @(90, "Neisseria", 0.0003247717504744928),
    @(92, "Porphyromonas", 0.00929406095537628),
    @(93, "Varibaculum", 0.006990870234154651),
    @(94, "Ezakiella", 0.002980975725860719),
    @(101, "Porphyromonas", 0.00929406095537628),
    @(102, "Neisseria", 0.0003247717504744928),
    @(105, "Porphyromonas", 0.00929406095537628),
    @(106, "Ezakiella", 0.002980975725860719),
    @(117, "Varibaculum", 0.006990870234154651),
    @(118, "Ezakiella", 0.002980975725860719),
    @(124, "Porphyromonas", 0.00929406095537628),
    @(125, "Haemophilus", 0.005704838249311408),
    @(126, "Neisseria", 0.0003247717504744928),
    @(173, "Porphyromonas", 0.00929406095537628),
    @(174, "Neisseria", 0.0003247717504744928),
    @(221, "Porphyromonas", 0.00929406095537628),
    @(222, "Neisseria", 0.0003247717504744928),
    @(225, "Porphyromonas", 0.00929406095537628),
    @(226, "Ezakiella", 0.002980975725860719),
    @(245, "Porphyromonas", 0.00929406095537628),
    @(246, "Neisseria", 0.0003247717504744928),
    @(248, "Varibaculum", 0.006990870234154651),
    @(249, "Ezakiella", 0.002980975725860719),
    @(280, "Porphyromonas", 0.00929406095537628),
    @(281, "Haemophilus", 0.005704838249311408),
    @(282, "Neisseria", 0.0003247717504744928),
    @(285, "Porphyromonas", 0.00929406095537628),
    @(286, "Ezakiella", 0.002980975725860719),
    @(292, "Porphyromonas", 0.00929406095537628),
    @(293, "Haemophilus", 0.005704838249311408),
    @(294, "Neisseria", 0.0003247717504744928),
    @(297, "Porphyromonas", 0.00929406095537628),
    @(298, "Ezakiella", 0.002980975725860719),
    @(309, "Varibaculum", 0.006990870234154651),
    @(310, "Ezakiella", 0.002980975725860719),
    @(364, "Porphyromonas", 0.00929406095537628),
    @(365, "Haemophilus", 0.005704838249311408),
    @(366, "Neisseria", 0.0003247717504744928),
    @(369, "Porphyromonas", 0.00929406095537628),
    @(370, "Ezakiella", 0.002980975725860719),
    @(377, "Porphyromonas", 0.00929406095537628),
    @(378, "Neisseria", 0.0003247717504744928),
    @(388, "Porphyromonas", 0.00929406095537628),
    @(389, "Haemophilus", 0.005704838249311408),
    @(390, "Neisseria", 0.0003247717504744928),
    @(393, "Porphyromonas", 0.00929406095537628),
    @(394, "Ezakiella", 0.002980975725860719),
    @(413, "Porphyromonas", 0.00929406095537628),
    @(414, "Neisseria", 0.0003247717504744928),
    @(429, "Varibaculum", 0.006990870234154651),
    @(430, "Ezakiella", 0.002980975725860719),
    @(448, "Porphyromonas", 0.00929406095537628),
    @(449, "Haemophilus", 0.005704838249311408),
    @(452, "Varibaculum", 0.006990870234154651),
    @(453, "Ezakiella", 0.002980975725860719),
    @(460, "Porphyromonas", 0.00929406095537628),
    @(461, "Haemophilus", 0.005704838249311408),
    @(477, "Varibaculum", 0.006990870234154651),
    @(478, "Ezakiella", 0.002980975725860719),
    @(497, "Porphyromonas", 0.00929406095537628),
    @(498, "Neisseria", 0.0003247717504744928),
    @(500, "Porphyromonas", 0.00929406095537628),
    @(501, "Varibaculum", 0.006990870234154651),
    @(502, "Ezakiella", 0.002980975725860719),
    @(525, "Varibaculum", 0.006990870234154651),
    @(526, "Ezakiella", 0.002980975725860719),
    @(533, "Porphyromonas", 0.00929406095537628),
    @(534, "Neisseria", 0.0003247717504744928),
    @(536, "Porphyromonas", 0.00929406095537628),
    @(537, "Varibaculum", 0.006990870234154651),
    @(538, "Ezakiella", 0.002980975725860719),
    @(544, "Porphyromonas", 0.00929406095537628),
    @(545, "Haemophilus", 0.005704838249311408),
    @(546, "Neisseria", 0.0003247717504744928),
    @(548, "Porphyromonas", 0.00929406095537628),
    @(549, "Varibaculum", 0.006990870234154651),
    @(550, "Ezakiella", 0.002980975725860719),
    @(556, "Porphyromonas", 0.00929406095537628),
    @(557, "Haemophilus", 0.005704838249311408),
    @(558, "Neisseria", 0.0003247717504744928),
    @(560, "Porphyromonas", 0.00929406095537628),
    @(561, "Varibaculum", 0.006990870234154651),
    @(562, "Ezakiella", 0.002980975725860719),
    @(568, "Porphyromonas", 0.00929406095537628),
    @(569, "Haemophilus", 0.005704838249311408),
    @(570, "Neisseria", 0.0003247717504744928),
    @(572, "Porphyromonas", 0.00929406095537628),
    @(573, "Varibaculum", 0.006990870234154651),
    @(574, "Ezakiella", 0.002980975725860719),
    @(581, "Porphyromonas", 0.00929406095537628),
    @(582, "Neisseria", 0.0003247717504744928),
    @(584, "Porphyromonas", 0.00929406095537628),
    @(585, "Varibaculum", 0.006990870234154651),
    @(586, "Ezakiella", 0.002980975725860719),
    @(592, "Porphyromonas", 0.00929406095537628),
    @(593, "Haemophilus", 0.005704838249311408),
    @(594, "Neisseria", 0.0003247717504744928),
    @(596, "Porphyromonas", 0.00929406095537628),
    @(597, "Varibaculum", 0.006990870234154651),
    @(598, "Ezakiella", 0.002980975725860719),
    @(604, "Porphyromonas", 0.00929406095537628),
    @(605, "Haemophilus", 0.005704838249311408),
    @(606, "Neisseria", 0.0003247717504744928),
    @(608, "Porphyromonas", 0.00929406095537628),
    @(609, "Varibaculum", 0.006990870234154651),
    @(610, "Ezakiella", 0.002980975725860719),
    @(628, "Porphyromonas", 0.00929406095537628),
    @(629, "Haemophilus", 0.005704838249311408),
    @(630, "Neisseria", 0.0003247717504744928),
    @(632, "Porphyromonas", 0.00929406095537628),
    @(633, "Varibaculum", 0.006990870234154651),
    @(634, "Ezakiella", 0.002980975725860719),
    @(645, "Varibaculum", 0.006990870234154651),
    @(646, "Ezakiella", 0.002980975725860719),
    @(652, "Porphyromonas", 0.00929406095537628),
    @(653, "Haemophilus", 0.005704838249311408),
    @(654, "Neisseria", 0.0003247717504744928),
    @(656, "Porphyromonas", 0.00929406095537628),
    @(657, "Varibaculum", 0.006990870234154651),
    @(658, "Ezakiella", 0.002980975725860719),
    @(665, "Porphyromonas", 0.00929406095537628),
    @(666, "Neisseria", 0.0003247717504744928),
    @(669, "Porphyromonas", 0.00929406095537628),
    @(670, "Ezakiella", 0.002980975725860719),
    @(705, "Varibaculum", 0.006990870234154651),
    @(706, "Ezakiella", 0.002980975725860719),
    @(712, "Porphyromonas", 0.00929406095537628),
    @(713, "Haemophilus", 0.005704838249311408),
    @(714, "Neisseria", 0.0003247717504744928),
    @(717, "Porphyromonas", 0.00929406095537628),
    @(718, "Ezakiella", 0.002980975725860719),
    @(737, "Porphyromonas", 0.00929406095537628),
    @(738, "Neisseria", 0.0003247717504744928),
    @(740, "Porphyromonas", 0.00929406095537628),
    @(741, "Varibaculum", 0.006990870234154651),
    @(742, "Ezakiella", 0.002980975725860719),
    @(748, "Porphyromonas", 0.00929406095537628),
    @(749, "Haemophilus", 0.005704838249311408),
    @(750, "Neisseria", 0.0003247717504744928),
    @(752, "Porphyromonas", 0.00929406095537628),
    @(753, "Varibaculum", 0.006990870234154651),
    @(754, "Ezakiella", 0.002980975725860719),
    @(760, "Porphyromonas", 0.00929406095537628),
    @(761, "Haemophilus", 0.005704838249311408),
    @(762, "Neisseria", 0.0003247717504744928),
    @(764, "Porphyromonas", 0.00929406095537628),
    @(765, "Varibaculum", 0.006990870234154651),
    @(766, "Ezakiella", 0.002980975725860719),
    @(789, "Varibaculum", 0.006990870234154651),
    @(790, "Ezakiella", 0.002980975725860719),
    @(808, "Porphyromonas", 0.00929406095537628),
    @(809, "Haemophilus", 0.005704838249311408),
    @(810, "Neisseria", 0.0003247717504744928),
    @(812, "Porphyromonas", 0.00929406095537628),
    @(813, "Varibaculum", 0.006990870234154651),
    @(814, "Ezakiella", 0.002980975725860719),
    @(821, "Porphyromonas", 0.00929406095537628),
    @(822, "Neisseria", 0.0003247717504744928),
    @(824, "Porphyromonas", 0.00929406095537628),
    @(825, "Varibaculum", 0.006990870234154651),
    @(826, "Ezakiella", 0.002980975725860719),
    @(856, "Porphyromonas", 0.00929406095537628),
    @(857, "Haemophilus", 0.005704838249311408),
    @(858, "Neisseria", 0.0003247717504744928),
    @(860, "Porphyromonas", 0.00929406095537628),
    @(861, "Varibaculum", 0.006990870234154651),
    @(862, "Ezakiella", 0.002980975725860719),
    @(868, "Porphyromonas", 0.00929406095537628),
    @(869, "Haemophilus", 0.005704838249311408),
    @(870, "Neisseria", 0.0003247717504744928),
    @(872, "Porphyromonas", 0.00929406095537628),
    @(873, "Varibaculum", 0.006990870234154651),
    @(874, "Ezakiella", 0.002980975725860719),
    @(885, "Varibaculum", 0.006990870234154651),
    @(886, "Ezakiella", 0.002980975725860719),
    @(904, "Porphyromonas", 0.00929406095537628),
    @(905, "Haemophilus", 0.005704838249311408),
    @(906, "Neisseria", 0.0003247717504744928),
    @(908, "Porphyromonas", 0.00929406095537628),
    @(909, "Varibaculum", 0.006990870234154651),
    @(910, "Ezakiella", 0.002980975725860719),
    @(928, "Porphyromonas", 0.00929406095537628),
    @(929, "Haemophilus", 0.005704838249311408),
    @(930, "Neisseria", 0.0003247717504744928),
    @(932, "Porphyromonas", 0.00929406095537628),
    @(933, "Varibaculum", 0.006990870234154651),
    @(934, "Ezakiella", 0.002980975725860719),
    @(941, "Porphyromonas", 0.00929406095537628),
    @(942, "Neisseria", 0.0003247717504744928),
    @(945, "Porphyromonas", 0.00929406095537628),
    @(946, "Ezakiella", 0.002980975725860719),
    @(964, "Porphyromonas", 0.00929406095537628),
    @(965, "Haemophilus", 0.005704838249311408),
    @(966, "Neisseria", 0.0003247717504744928),
    @(968, "Porphyromonas", 0.00929406095537628),
    @(969, "Varibaculum", 0.006990870234154651),
    @(970, "Ezakiella", 0.002980975725860719),
    @(993, "Varibaculum", 0.006990870234154651),
    @(994, "Ezakiella", 0.002980975725860719),
    @(1012, "Porphyromonas", 0.00929406095537628),
    @(1013, "Haemophilus", 0.005704838249311408),
    @(1014, "Neisseria", 0.0003247717504744928),
    @(1017, "Porphyromonas", 0.00929406095537628),
    @(1018, "Ezakiella", 0.002980975725860719),
    @(1029, "Varibaculum", 0.006990870234154651),
    @(1030, "Ezakiella", 0.002980975725860719),
    @(1049, "Porphyromonas", 0.00929406095537628),
    @(1050, "Neisseria", 0.0003247717504744928),
    @(1061, "Porphyromonas", 0.00929406095537628),
    @(1062, "Neisseria", 0.0003247717504744928),
    @(1065, "Porphyromonas", 0.00929406095537628),
    @(1066, "Ezakiella", 0.002980975725860719),
    @(1072, "Porphyromonas", 0.00929406095537628),
    @(1073, "Haemophilus", 0.005704838249311408),
    @(1074, "Neisseria", 0.0003247717504744928),
    @(1076, "Porphyromonas", 0.00929406095537628),
    @(1077, "Varibaculum", 0.006990870234154651),
    @(1078, "Ezakiella", 0.002980975725860719),
    @(1084, "Porphyromonas", 0.00929406095537628),
    @(1085, "Haemophilus", 0.005704838249311408),
    @(1086, "Neisseria", 0.0003247717504744928),
    @(1088, "Porphyromonas", 0.00929406095537628),
    @(1089, "Varibaculum", 0.006990870234154651),
    @(1090, "Ezakiella", 0.002980975725860719)
)

foreach ($item in $changes) {
    $row = $item[0]
    $species = $item[1]
    $eVal = $item[2]
    $ws.Cells.Item($row, 3).Value = $species
    $ws.Cells.Item($row, 5).Value = $eVal
}

Write-Host ("Updated " + $changes.Count + " rows")
